# Applies the "Add files via upload" commit:
#  1. Renames the shared string "Orta" -> "Orta Çember" (affects AK185:AK188,
#     which are the only cells currently holding the plain "Orta" text).
#  2. Appends four new game rows (220-223) to Sheet1 with full data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) "Orta" -> "Orta Çember"
# ---------------------------------------------------------------------------
$AK = 37   # column AK
foreach ($r in 185..188) {
    $ws.Cells.Item($r, $AK).Value = "Orta Çember"
}

# ---------------------------------------------------------------------------
# 2) Append new rows 220-223
# ---------------------------------------------------------------------------
$newRows = @(220, 221, 222, 223)
$srcRow = 211   # a fully-populated row (A:D, F:AK) to source cell formats from

foreach ($r in $newRows) {
    # Copy number/date/time formats from an existing row so the new cells
    # land on the same style slots (s="1"/"4"/"5") instead of minting brand
    # new styles. Column E is intentionally excluded - it is always blank
    # in this sheet.
    $ws.Range("A$srcRow`:D$srcRow").Copy()
    $ws.Range("A$r`:D$r").PasteSpecial(-4122)
    $ws.Range("F$srcRow`:AK$srcRow").Copy()
    $ws.Range("F$r`:AK$r").PasteSpecial(-4122)
}

$rows = @{
    220 = @{
        A=53.0; B=44846.0; C="Ecem"; D=1.0; F="Batuhan"; G=9.0; H="Turuncu"; I=4.0; J=1.0; K=7.0; L=1.0;
        M=0.0; N=0.0; O=0.0; P=2.0; Q="Saman, Odun"; R=1.0; S=50.0; T=0.0; U=0.5694444444444444;
        V=0.0; W=0.0; X=0.0; Y=0.0; Z=0.0; AA=0.0; AB=0.0; AC=0.0; AD=0.0; AE=1.0; AF=1.0; AG=0.0; AH=0.0; AI=0.0;
        AJ="3, 4, 6, 8, 4"; AK="Orta Çember"
    };
    221 = @{
        A=53.0; B=44846.0; C="Batuhan"; D=2.0; F="Batuhan"; G=10.0; H="Mavi"; I=1.0; J=5.0; K=9.0; L=5.0;
        M=0.0; N=0.0; O=1.0; P=2.0; Q="Odun, Saman, Koyun"; R=0.0; S=50.0; T=0.0; U=0.5694444444444444;
        V=0.0; W=1.0; X=0.0; Y=0.0; Z=0.0; AA=0.0; AB=0.0; AC=0.0; AD=1.0; AE=0.0; AF=0.0; AG=0.0; AH=0.0; AI=1.0;
        AJ="6, 9, 10, 2, 5, 9"; AK="Orta Çember"
    };
    222 = @{
        A=53.0; B=44846.0; C="Çağatay"; D=3.0; F="Batuhan"; G=8.0; H="Kırmızı"; I=0.0; J=5.0; K=11.0; L=2.0;
        M=1.0; N=0.0; O=0.0; P=1.0; Q="Tuğla, Koyun, Koyun"; R=0.0; S=50.0; T=0.0; U=0.5694444444444444;
        V=1.0; W=0.0; X=0.0; Y=0.0; Z=0.0; AA=0.0; AB=0.0; AC=0.0; AD=0.0; AE=0.0; AF=0.0; AG=1.0; AH=0.0; AI=0.0;
        AJ="3, 8, 10, 5, 6, 11"; AK="Orta Çember"
    };
    223 = @{
        A=53.0; B=44846.0; C="Alperen"; D=4.0; F="Batuhan"; G=8.0; H="Beyaz"; I=4.0; J=0.0; K=6.0; L=0.0;
        M=0.0; N=1.0; O=0.0; P=2.0; Q="Tuğla, Tuğla, Koyun"; R=0.0; S=50.0; T=0.0; U=0.5694444444444444;
        V=0.0; W=0.0; X=0.0; Y=0.0; Z=0.0; AA=0.0; AB=0.0; AC=0.0; AD=0.0; AE=0.0; AF=2.0; AG=0.0; AH=0.0; AI=0.0;
        AJ="5, 8, 10, 4, 9, 11"; AK="Orta Çember"
    }
}

$colNum = @{
    "A"=1; "B"=2; "C"=3; "D"=4; "F"=6; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11; "L"=12; "M"=13; "N"=14; "O"=15;
    "P"=16; "Q"=17; "R"=18; "S"=19; "T"=20; "U"=21; "V"=22; "W"=23; "X"=24; "Y"=25; "Z"=26; "AA"=27; "AB"=28;
    "AC"=29; "AD"=30; "AE"=31; "AF"=32; "AG"=33; "AH"=34; "AI"=35; "AJ"=36; "AK"=37
}

foreach ($r in $newRows) {
    $rowData = $rows[$r]
    foreach ($col in $rowData.Keys) {
        $c = $colNum[$col]
        $ws.Cells.Item($r, $c).Value = $rowData[$col]
    }
}
